# Delete row 364 ("「開けゴマ」..." post) from Sheet1.
# This causes all subsequent rows to shift up by one and the sheet's
# dimension to shrink from A1:C567 to A1:C566, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(364).Delete()

$wb.Save()
